# array practice questions implementation
$wb = $excel.ActiveWorkbook

# --- Rename 3rd sheet "Practice Qns" -> "Practice Questions" ---
$practice = $wb.Worksheets.Item("Practice Qns")
$practice.Name = "Practice Questions"

$pythonCode = $wb.Worksheets.Item("pythonCode")

# ============================================================
# Sheet "pythonCode" (physically sheet2.xml)
# ============================================================
# Row1: A1 header changes from "TestId" to "pythonCode"
$pythonCode.Range("A1").Value = "pythonCode"
$pythonCode.Range("A1").Style = "Normal 2"

# Row2: clear out old TC001/hello values (headers row stays the same style)
$pythonCode.Range("A2").Value = ""
$pythonCode.Range("A2").Style = "Normal 2"
$pythonCode.Range("C2").Value = ""

# Row3: becomes the "hello" test row
$pythonCode.Range("A3").Value = "hello"
$pythonCode.Range("A3").Style = "Normal 2"
$pythonCode.Range("C3").Value = "NameError: name 'hello' is not defined on line 1"

# Row4: becomes the print("hello") test row
$pythonCode.Range("A4").Value = 'print("hello")'
$pythonCode.Range("C4").Value = "hello"

# Rows 5-11: clear leftover "submission success"/"2"/etc values
$pythonCode.Range("C5").Value = ""
$pythonCode.Range("C6").Value = ""
$pythonCode.Range("C7").Value = ""
$pythonCode.Range("C8").Value = ""
$pythonCode.Range("C9").Value = ""
$pythonCode.Range("C10").Value = ""
$pythonCode.Range("C11").Value = ""

# New font for the "Result" column error cells (Consolas, 10pt, black)
$f = $pythonCode.Range("C3").Font
$f.Name = "Consolas"
$f.Size = 10
$f.Color = 0

# View: show column C, select C1:C4
$pythonCode.Range("C1:C4").Select()

# ============================================================
# Sheet "Practice Questions" (physically sheet3.xml)
# ============================================================
$practice.Range("C1").Value = "Result"
$practice.Range("C1").Style = "Normal 2"
$practice.Range("C2").Value = ""
$practice.Range("C2").Style = "Normal 2"
$practice.Range("C3").Value = "NameError: name 'hello' is not defined on line 1"
$practice.Range("C4").Value = "hello"
$practice.Range("C4").Style = "Normal 2"
$practice.Range("B4").Value = 'print("hello")'

$f2 = $practice.Range("C3").Font
$f2.Name = "Consolas"
$f2.Size = 10
$f2.Color = 0

$practice.Range("B3").Select()
